# Update crypto price/volume data per latest scrape (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.194.10'
$ws.Cells.Item(2, 5).Value = '  -1.50%  '
$ws.Cells.Item(3, 4).Value = '1.805.15'
$ws.Cells.Item(3, 5).Value = '  +0.37%  '
$c = $ws.Cells.Item(4, 4)
$c.Value = "'1.0000"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$c = $ws.Cells.Item(5, 4)
$c.Value = "'316.78"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.03%  '
$c = $ws.Cells.Item(6, 4)
$c.Value = "'0.9998"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.11%  '
$c = $ws.Cells.Item(7, 4)
$c.Value = "'0.5375"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +1.20%  '
$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.3785"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +0.65%  '
$c = $ws.Cells.Item(9, 4)
$c.Value = "'0.07485"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -0.63%  '
$c = $ws.Cells.Item(10, 4)
$c.Value = "'42.01"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -1.21%  '
$c = $ws.Cells.Item(11, 4)
$c.Value = "'1.099"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -2.12%  '
$c = $ws.Cells.Item(12, 4)
$c.Value = "'0.9998"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.10%  '
$c = $ws.Cells.Item(13, 4)
$c.Value = "'6.210"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.13%  '
$c = $ws.Cells.Item(14, 4)
$c.Value = "'20.54"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -3.12%  '
$c = $ws.Cells.Item(15, 4)
$c.Value = "'7.392"
$c.Style = "Normal"
$ws.Cells.Item(16, 4).Value = '1.803.84'
$ws.Cells.Item(16, 5).Value = '  +0.68%  '
$c = $ws.Cells.Item(17, 4)
$c.Value = "'89.95"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.67%  '
$c = $ws.Cells.Item(18, 4)
$c.Value = "'0.00001066"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 4)
$c.Value = "'0.06499"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.71%  '
$c = $ws.Cells.Item(20, 4)
$c.Value = "'17.43"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.63%  '
$c = $ws.Cells.Item(21, 4)
$c.Value = "'0.9994"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.06%  '
$c = $ws.Cells.Item(22, 4)
$c.Value = "'5.934"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.12%  '
$ws.Cells.Item(23, 4).Value = '28.213.61'
$ws.Cells.Item(23, 5).Value = '  -1.46%  '
$c = $ws.Cells.Item(24, 4)
$c.Value = "'11.22"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.04%  '
$c = $ws.Cells.Item(25, 4)
$c.Value = "'2.091"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.17%  '
$c = $ws.Cells.Item(26, 4)
$c.Value = "'155.98"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -3.09%  '
$ws.Cells.Item(27, 5).Value = '  -0.13%  '
$ws.Cells.Item(28, 4).Value = '2.011.43'
$ws.Cells.Item(28, 5).Value = '  +0.68%  '
$c = $ws.Cells.Item(29, 4)
$c.Value = "'2.341"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -2.44%  '
$c = $ws.Cells.Item(30, 4)
$c.Value = "'122.26"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -1.12%  '
$c = $ws.Cells.Item(31, 4)
$c.Value = "'1.129"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -0.37%  '
$c = $ws.Cells.Item(32, 4)
$c.Value = "'0.1117"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +9.30%  '
$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Cells.Item(33, 4)
$c.Value = "'5.619"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.85%  '
$ws.Cells.Item(34, 2).Value = 'HuobiToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Cells.Item(34, 4)
$c.Value = "'3.672"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.20%  '
$c = $ws.Cells.Item(35, 4)
$c.Value = "'0.06994"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +6.63%  '
$c = $ws.Cells.Item(36, 4)
$c.Value = "'0.2232"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -3.08%  '
$ws.Cells.Item(37, 5).Value = '  -0.91%  '
$c = $ws.Cells.Item(38, 4)
$c.Value = "'5.106"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +0.52%  '
$c = $ws.Cells.Item(39, 4)
$c.Value = "'8.473"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -3.88%  '
$c = $ws.Cells.Item(40, 4)
$c.Value = "'11.20"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -2.68%  '
$c = $ws.Cells.Item(41, 4)
$c.Value = "'0.6184"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -2.14%  '
$c = $ws.Cells.Item(42, 4)
$c.Value = "'1.429"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +2.42%  '
$c = $ws.Cells.Item(43, 4)
$c.Value = "'1.177"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -1.85%  '
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(44, 4)
$c.Value = "'13.49"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -0.12%  '
$ws.Cells.Item(45, 2).Value = 'PancakeSwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Cells.Item(45, 4)
$c.Value = "'3.686"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.54%  '
$ws.Cells.Item(46, 2).Value = 'Decentraland'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Cells.Item(46, 4)
$c.Value = "'0.5782"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -2.44%  '
$ws.Cells.Item(47, 2).Value = 'Quant'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Cells.Item(47, 4)
$c.Value = "'125.66"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -0.13%  '
$ws.Cells.Item(48, 2).Value = 'EOS'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$c = $ws.Cells.Item(48, 4)
$c.Value = "'1.190"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +1.56%  '
$ws.Cells.Item(49, 2).Value = 'NEARProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Cells.Item(49, 4)
$c.Value = "'1.933"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -2.47%  '
$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Cells.Item(50, 4)
$c.Value = "'0.06825"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -1.60%  '
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Cells.Item(51, 4)
$c.Value = "'72.04"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -1.45%  '
